$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Intercept)
$ws.Range("B2").Value = 2.7573
$ws.Range("C2").Value = 0.08890000000000001
$ws.Range("D2").Value = 31.0211
$ws.Range("F2").Value = 2.5806
$ws.Range("G2").Value = 2.934

# Row 3 (Effort Expectancy)
$ws.Range("B3").Value = 0.5178
$ws.Range("C3").Value = 0.09039999999999999
$ws.Range("D3").Value = 5.7312
$ws.Range("F3").Value = 0.3382
$ws.Range("G3").Value = 0.6975

# Row 4 (AI Use Frequency)
$ws.Range("B4").Value = 0.0367
$ws.Range("C4").Value = 0.0723
$ws.Range("D4").Value = 0.5077
$ws.Range("E4").Value = 0.613
$ws.Range("F4").Value = -0.107
$ws.Range("G4").Value = 0.1804

# Row 5 (Interaction)
$ws.Range("B5").Value = 0.0495
$ws.Range("C5").Value = 0.0602
$ws.Range("D5").Value = 0.8228
$ws.Range("E5").Value = 0.4129
$ws.Range("F5").Value = -0.0701
$ws.Range("G5").Value = 0.1691
